$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 3 with the latest prediction results (ticker appears with results)
$ws.Range("A3").Value = 42632.880844907406
$ws.Range("B3").Value = 8
$ws.Range("C3").Value = "Buy"
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = "Random"
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0.56000000000000005
$ws.Range("S3").Value = 0.10150000000000001
$ws.Range("T3").Value = -0.93
$ws.Range("U3").Value = 2.3199999999999998
$ws.Range("V3").Value = "N/A"
$ws.Range("W3").Value = 0

# Match number formatting used on row 2 (date format on column A, percent format on column S)
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)

$ws.Range("S2").Copy()
$ws.Range("S3").PasteSpecial(-4122)

$excel.CutCopyMode = 0
